$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.737.11"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3
$ws.Range("D3").Value = "2.802.06"
$ws.Range("E3").Value = "  +0.80%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "354.03"
$ws.Range("E5").Value = "  -0.95%  "

# Row 6
$ws.Range("D6").Value = "111.45"
$ws.Range("E6").Value = "  +2.17%  "

# Row 7
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").Value = "  +0.66%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  +7.47%  "

# Row 10
$ws.Range("D10").Value = "40.30"
$ws.Range("E10").Value = "  +1.47%  "

# Row 11
$ws.Range("E11").Value = "  -2.90%  "

# Row 12
$ws.Range("E12").Value = "  -0.52%  "

# Row 13
$ws.Range("D13").Value = "20.02"
$ws.Range("E13").Value = "  +1.27%  "

# Row 14
$ws.Range("E14").Value = "  +2.19%  "

# Row 15
$ws.Range("D15").Value = "3.240.97"
$ws.Range("E15").Value = "  +0.81%  "

# Row 16
$ws.Range("D16").Value = "2.804.12"
$ws.Range("E16").Value = "  +0.79%  "

# Row 17
$ws.Range("E17").Value = "  +2.44%  "

# Row 18
$ws.Range("D18").Value = "51.672.21"
$ws.Range("E18").Value = "  +0.23%  "

# Row 19
$ws.Range("D19").Value = "7.62"
$ws.Range("E19").Value = "  +0.37%  "

# Row 20
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "13.92"
$ws.Range("E20").Value = "  +5.47%  "

# Row 21
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "3.19"
$ws.Range("E21").Value = "  +3.34%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  +0.40%  "

# Row 23
$ws.Range("D23").Value = "70.33"
$ws.Range("E23").Value = "  +0.29%  "

# Row 24
$ws.Range("D24").Value = "267.57"
$ws.Range("E24").Value = "  -0.24%  "

# Row 25
$ws.Range("D25").Value = "2.77"
$ws.Range("E25").Value = "  +0.88%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.01%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "26.16"
$ws.Range("E27").Value = "  -0.63%  "

# Row 28
$ws.Range("E28").Value = "  -3.10%  "

# Row 29
$ws.Range("D29").Value = "39.26"
$ws.Range("E29").Value = "  +11.52%  "

# Row 30
$ws.Range("D30").Value = "10.35"
$ws.Range("E30").Value = "  +1.72%  "

# Row 31
$ws.Range("E31").Value = "  +2.70%  "

# Row 32
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "52.46"
$ws.Range("E32").Value = "  +0.91%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "6.14"
$ws.Range("E33").Value = "  -0.94%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.0889"
$ws.Range("E34").Value = "  +6.16%  "

# Row 35
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "5.62"
$ws.Range("E35").Value = "  +8.15%  "

# Row 36
$ws.Range("E36").Value = "  +1.10%  "

# Row 37
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.06%  "

# Row 38
$ws.Range("D38").Value = "18.91"
$ws.Range("E38").Value = "  +0.98%  "

# Row 39
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "3.17"
$ws.Range("E39").Value = "  +1.20%  "

# Row 40
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "2.01"
$ws.Range("E40").Value = "  +2.90%  "

# Row 41
$ws.Range("E41").Value = "  +0.90%  "

# Row 42
$ws.Range("D42").Value = "2.50"
$ws.Range("E42").Value = "  -0.22%  "

# Row 43
$ws.Range("E43").Value = "  +2.01%  "

# Row 44
$ws.Range("D44").Value = "121.11"
$ws.Range("E44").Value = "  +1.55%  "

# Row 45
$ws.Range("D45").Value = "21.81"
$ws.Range("E45").Value = "  +0.15%  "

# Row 46
$ws.Range("D46").Value = "2.47"
$ws.Range("E46").Value = "  +6.76%  "

# Row 47
$ws.Range("E47").Value = "  +5.12%  "

# Row 48
$ws.Range("D48").Value = "2.107.09"
$ws.Range("E48").Value = "  +0.94%  "

# Row 49
$ws.Range("D49").Value = "0.954"
$ws.Range("E49").Value = "  +2.03%  "

# Row 50
$ws.Range("D50").Value = "5.47"
$ws.Range("E50").Value = "  -1.32%  "

# Row 51
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "0.219"
$ws.Range("E51").Value = "  +16.63%  "

